# LOM3261.docx edit: rotate the Objetivos / Docente(s) / Programa resumido /
# Programa / Avaliacao / Bibliografia text blocks into their new positions,
# as described by the supplied unified diff. Paragraph structure, styles and
# run formatting (bold/italic) stay the same; only the w:t contents (and, in
# a few spots, how a run's text is split across w:br line breaks) change.
#
# Every Find/Replace below is scoped to a specific paragraph's Range so that
# the operations are independent of execution order and cannot bleed into
# neighboring runs/paragraphs that happen to share text.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Replace failed in paragraph $paraIndex for text: $oldText"
    }
}

# --- Paragraph 6: Objetivos (PT body) ---
Replace-InParagraph 6 `
    "Introduzir o uso e prática de métodos numéricos para a solução de problemas matemáticos aplicados à Física e Engenharia. O estudante estará capacitado a descrever matematicamente e resolver numericamente problemas com o auxílio de algoritmos computacionais." `
    "Representação computacional de números em ponto flutuante; Zeros de funções; Sistemas de equações lineares; Método dos Mínimos Quadrados; Interpolação; Integração numérica; equações diferenciais ordinárias."

# --- Paragraph 7: Objetivos (EN, italic body) ---
Replace-InParagraph 7 `
    "Introduce the use and practice of numerical methods for solving mathematical problems applied to Physics and Engineering. The student will be able to describe mathematically and numerically solve problems with the aid of computational algorithms." `
    "Computational representation of floating point numbers; Root finding; Systems of linear equations; Least squares; Interpolation; Numerical integration; ordinary differential equations."

# --- Paragraph 9: Docente(s) Responsável(eis) list (two runs) ---
Replace-InParagraph 9 `
    "3480026 - João Paulo Pascon" `
    "Introduzir o uso e prática de métodos numéricos para a solução de problemas matemáticos aplicados à Física e Engenharia. O estudante estará capacitado a descrever matematicamente e resolver numericamente problemas com o auxílio de algoritmos computacionais."

Replace-InParagraph 9 `
    "1176388 - Luiz Tadeu Fernandes Eleno" `
    "• Números em ponto flutuante: representação e precisão. • Raízes de funções: método da bissecção; método da falsa posição; método de Newton-Raphson; • Solução de sistemas de equações lineares: pivotamento e escalonamento; método de Gauss. • Método dos mínimos quadrados: ajuste de funções lineares nos parâmetros ajustáveis; ajuste de funções linearizáveis; ajuste de funções não-lineares usando a biblioteca scipy.optimize • Interpolação: método de Lagrange; método de Newton; • Integração numérica: regra dos trapézios; regra de Simpson; métodos avançados implementados na biblioteca scipy.integrate.  • Solução de equações diferenciais ordinárias: método de Euler; método de Runge-Kutta; métodos mais avançados da biblioteca scipy.integrate."

# --- Paragraph 11: Programa resumido (PT body) ---
Replace-InParagraph 11 `
    "Representação computacional de números em ponto flutuante; Zeros de funções; Sistemas de equações lineares; Método dos Mínimos Quadrados; Interpolação; Integração numérica; equações diferenciais ordinárias." `
    "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados."

# --- Paragraph 12: Programa resumido (EN, italic body) ---
Replace-InParagraph 12 `
    "Computational representation of floating point numbers; Root finding; Systems of linear equations; Least squares; Interpolation; Numerical integration; ordinary differential equations." `
    "Introduce the use and practice of numerical methods for solving mathematical problems applied to Physics and Engineering. The student will be able to describe mathematically and numerically solve problems with the aid of computational algorithms."

# --- Paragraph 14: Programa (PT body) ---
Replace-InParagraph 14 `
    "• Números em ponto flutuante: representação e precisão. • Raízes de funções: método da bissecção; método da falsa posição; método de Newton-Raphson; • Solução de sistemas de equações lineares: pivotamento e escalonamento; método de Gauss. • Método dos mínimos quadrados: ajuste de funções lineares nos parâmetros ajustáveis; ajuste de funções linearizáveis; ajuste de funções não-lineares usando a biblioteca scipy.optimize • Interpolação: método de Lagrange; método de Newton; • Integração numérica: regra dos trapézios; regra de Simpson; métodos avançados implementados na biblioteca scipy.integrate.  • Solução de equações diferenciais ordinárias: método de Euler; método de Runge-Kutta; métodos mais avançados da biblioteca scipy.integrate." `
    "Média aritmética de trabalhos propostos ao longo do curso (30%) e duas avaliações individuais (70%)."

# --- Paragraph 17: Avaliação list (Método: / Critério: / Norma de recuperação:) ---
Replace-InParagraph 17 `
    "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados." `
    "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

Replace-InParagraph 17 `
    "Média aritmética de trabalhos propostos ao longo do curso (30%) e duas avaliações individuais (70%)." `
    ("Cunha, M. C. C., Métodos Numéricos. Editora Unicamp, 1993.^l" + `
     "Sperandio, D., Mendes, J. T., Monken e Silva, L. H. Cálculo Numérico. Pearson, 2003 ^l" + `
     "LANGTANGEN, Hans Petter. A Primer on scientific programming with Python, 2a ed. New York: Springer, 2011. ^l" + `
     "LANGTANGEN, Hans Petter. Python scripting for computational science, 5a ed. New York: Springer, 2016. ^l" + `
     "SCOPATZ, A.; HUFF, K. D. Effective computation in physics: field guide to research in Python. Sebastpol, CA: O’Reilly Media, 2015.")

Replace-InParagraph 17 `
    "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação" `
    "3480026 - João Paulo Pascon"

# --- Paragraph 19: Bibliografia body (five text+break runs collapse to one) ---
# The whole paragraph content spans multiple w:t/w:br pairs inside one run, so
# use a wildcard Find across the full span rather than a plain substring match.
$rng19 = $d.Paragraphs(19).Range
$f19 = $rng19.Find
$f19.Text = "Cunha, M. C. C.*Sebastpol, CA: O’Reilly Media, 2015."
$f19.MatchWildcards = $true
$f19.Replacement.Text = "1176388 - Luiz Tadeu Fernandes Eleno"
$ok19 = $f19.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)
if (-not $ok19) {
    throw "Replace failed in paragraph 19 (Bibliografia block)"
}

Write-Output "done"
